$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timing entry: training classified on all training samples, saving model results.
$ws.Range("A17").Value = "1590984652.9147887"

$ws.Range("A18").Value = "1590984770.601893"
$ws.Range("B18").Formula = "=A18-A17"
$ws.Range("C18").Formula = "=B18/60"
$ws.Range("D18").Formula = "=310000/1000*C18/60"
$ws.Range("D18").NumberFormat = "0.00"
$ws.Range("E18").Value = "on vast.ai machine with Pool() (having 16 cores)"
$ws.Range("F18").Value = "docker"
$ws.Range("H18").Value = "there is another machine with 48 cores"

# Highlight the new cost estimate and machine note with a yellow fill.
$ws.Range("D18:E18").Interior.Color = 65535

# Update the active selection to reflect where the author left off.
$ws.Range("I12").Select() | Out-Null
